$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 6 (row 14): C14/D14/E14 already had values; fill in the rest.
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 2

# Week 7 (row 16): C16 gates the row total formula, so set it last
# so the recalculated SUM() picks up all the other values in the row.
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("C16").Value = 2

# Week 8 (row 18): same approach, C18 last.
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.5
$ws.Range("F18").Value = 0.5
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("C18").Value = 2

# Week 9 (row 20): same approach, C20 last. (I20 stays blank.)
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("C20").Value = 0

# Update the active selection to reflect the last edited cell.
$ws.Range("I20").Select()
